# Update countries & provincias Spain
# Applies the daily COVID-19 data refresh to the "Pais" sheet:
#  - refreshes case/death/recovery counters for several countries
#  - Tanzania's updated numbers move it above Martinica & Guadalupe in the
#    (descending, by "Casos totales") ranking, so those two rows keep their
#    own figures but shift down one row
#  - bumps the "Datos actualizados ..." timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 740928
$ws.Range("C4").Value = 2136
$ws.Range("D4").Value = 68599
$ws.Range("E4").Value = 633245
$ws.Range("G4").Value = 70
$ws.Range("H4").Value = 39084

# Row 20
$ws.Range("B20").Value = 16402
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 2601
$ws.Range("E20").Value = 13263
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 538

# Row 22
$ws.Range("B22").Value = 14699
$ws.Range("C22").Value = 28
$ws.Range("E22").Value = 3755

# Row 39
$ws.Range("B39").Value = 6612
$ws.Range("C39").Value = 26
$ws.Range("E39").Value = 2311
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 71

# Row 47
$ws.Range("E47").Value = 4390
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 47

# Row 48
$ws.Range("B48").Value = 4680
$ws.Range("C48").Value = 345
$ws.Range("D48").Value = 363
$ws.Range("E48").Value = 4091
$ws.Range("G48").Value = 9
$ws.Range("H48").Value = 226

# Row 66
$ws.Range("B66").Value = 1771
$ws.Range("C66").Value = 11
$ws.Range("E66").Value = 471

# Row 88
$ws.Range("B88").Value = 767
$ws.Range("C88").Value = 6
$ws.Range("E88").Value = 676

# Row 115
$ws.Range("D115").Value = 193
$ws.Range("E115").Value = 99
$ws.Range("F115").Value = 10

# Rows 127-129: Tanzania's refreshed numbers (170) overtake Martinica (163)
# in the "Casos totales" ranking, so Tanzania now sits in row 127, and
# Martinica / Guadalupe shift down to rows 128 / 129 keeping their own data.
$ws.Range("A127").Value = "Tanzania"
$ws.Range("B127").Value = 170
$ws.Range("C127").Value = 23
$ws.Range("D127").Value = 11
$ws.Range("E127").Value = 152
$ws.Range("F127").Value = 4
$ws.Range("G127").Value = 2
$ws.Range("H127").Value = 7

$ws.Range("A128").Value = "Martinica"
$ws.Range("B128").Value = 163
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 73
$ws.Range("E128").Value = 78
$ws.Range("F128").Value = 11
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 12

$ws.Range("A129").Value = "Guadalupe"
$ws.Range("B129").Value = 148
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 73
$ws.Range("E129").Value = 67
$ws.Range("F129").Value = 13
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 8

# Row 136
$ws.Range("B136").Value = 121
$ws.Range("C136").Value = 1
$ws.Range("D136").Value = 39
$ws.Range("E136").Value = 82

# Row 159
$ws.Range("B159").Value = 52
$ws.Range("C159").Value = 17
$ws.Range("D159").Value = 17
$ws.Range("E159").Value = 36

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 17:22"
